$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("B2").Value = 1.0
$ws.Range("D2").Value = 22.0
$ws.Range("E2").Value = 23.0
$ws.Range("F2").Value = 986.52
$ws.Range("G2").Value = 23.62

# Add new row 3
$ws.Range("A3").Value = "Viagra"
$ws.Range("B3").Value = 1234.0
$ws.Range("C3").Value = "R"
$ws.Range("D3").Value = 23.0
$ws.Range("E3").Value = 22.0
$ws.Range("J3").Value = 1000.0
$ws.Range("K3").Value = 10.0
$ws.Range("L3").Value = 91.0
$ws.Range("M3").Value = 9.0
$ws.Range("N3").Value = 0.0
$ws.Range("O3").Value = 20.0
